# Workbook shows a Gantt-style planning sheet. A second "Réalisé" pass was
# marked (column C holds a merged "V" per finished task group, column D the
# date it was done) for the groups starting at rows 49, 54 and 62, plus the
# scroll position / active cell were updated as the user kept working down
# the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Column C formatting: reuse the centered alignment already used by the
# C43:C46 "V" group. Apply it across the whole touched span (C49:C79) so
# the per-row column span bookkeeping matches, then clear the cells that
# must remain completely untouched (blank separator rows, and rows whose
# only change is bookkeeping, not an actual new cell).
# ---------------------------------------------------------------------
$ws.Range("C44").Copy()
$ws.Range("C49:C79").PasteSpecial(-4122)

# Blank separator rows: must stay fully empty (no cell at all).
$ws.Range("C52").Clear()
$ws.Range("C60").Clear()
$ws.Range("C66").Clear()
$ws.Range("C75").Clear()

# Rows that already have other content but must NOT gain a new C cell.
$ws.Range("C53").Clear()
$ws.Range("C61").Clear()
$ws.Range("C67:C74").Clear()
$ws.Range("C76:C79").Clear()

# ---------------------------------------------------------------------
# "V" markers (shared string already used elsewhere in the sheet) on the
# first row of each newly finished task group.
# ---------------------------------------------------------------------
$ws.Range("C49").Value = "V"
$ws.Range("C54").Value = "V"
$ws.Range("C62").Value = "V"

# ---------------------------------------------------------------------
# Column D: completion dates for the newly finished rows.
# ---------------------------------------------------------------------
$ws.Range("D51").Value = 43917
$ws.Range("D51").NumberFormat = "d-mmm"

$ws.Range("D54:D59").Value = 43917
$ws.Range("D54:D59").NumberFormat = "d-mmm"

$ws.Range("D62:D65").Value = 43917
$ws.Range("D62:D65").NumberFormat = "d-mmm"

# ---------------------------------------------------------------------
# Merge each "V" cell across its task group, like the existing C43:C46.
# ---------------------------------------------------------------------
$ws.Range("C49:C51").Merge()
$ws.Range("C54:C59").Merge()
$ws.Range("C62:C65").Merge()

# ---------------------------------------------------------------------
# View state: the user scrolled further down and moved the selection.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K67").Select()
